$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" '29.387.33'
Set-TextValue "E2" '  +0.14%  '

Set-TextValue "D3" '1.884.41'
Set-TextValue "E3" '  +0.41%  '

Set-TextValue "E4" '  +0.12%  '

Set-TextValue "D5" '0.7126'
Set-TextValue "E5" '  +0.17%  '

Set-TextValue "D6" '242.40'
Set-TextValue "E6" '  +0.05%  '

Set-TextValue "E7" '  +0.08%  '

Set-TextValue "D8" '0.08034'
Set-TextValue "E8" '  +3.63%  '

Set-TextValue "D9" '0.3128'
Set-TextValue "E9" '  +0.57%  '

Set-TextValue "D10" '25.25'
Set-TextValue "E10" '  +0.78%  '

Set-TextValue "D11" '0.08358'
Set-TextValue "E11" '  -1.24%  '

Set-TextValue "D12" '1.897.60'
Set-TextValue "E12" '  -1.84%  '

Set-TextValue "D13" '0.7208'
Set-TextValue "E13" '  +1.28%  '

Set-TextValue "D14" '5.246'
Set-TextValue "E14" '  +0.68%  '

Set-TextValue "D15" '92.59'
Set-TextValue "E15" '  +1.29%  '

Set-TextValue "D16" '6.300'
Set-TextValue "E16" '  +4.93%  '

Set-TextValue "D17" '0.000008482'
Set-TextValue "E17" '  +2.45%  '

Set-TextValue "D18" '29.413.94'
Set-TextValue "E18" '  +0.20%  '

Set-TextValue "B19" 'WrappedliquidstakedEther2.0'
Set-TextValue "C19" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D19" '2.156.41'
Set-TextValue "E19" '  +1.14%  '

Set-TextValue "B20" 'BitcoinCash'
Set-TextValue "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D20" '241.28'
Set-TextValue "E20" '  -0.48%  '

Set-TextValue "B21" 'Avalanche'
Set-TextValue "C21" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D21" '13.26'
Set-TextValue "E21" '  +0.17%  '

Set-TextValue "E22" '  +0.10%  '

Set-TextValue "D23" '7.871'
Set-TextValue "E23" '  +0.25%  '

Set-TextValue "E24" '  +0.09%  '

Set-TextValue "D25" '0.1587'
Set-TextValue "E25" '  -1.21%  '

Set-TextValue "D26" '163.74'
Set-TextValue "E26" '  +0.59%  '

Set-TextValue "D27" '9.076'
Set-TextValue "E27" '  +0.63%  '

Set-TextValue "D29" '1.506'
Set-TextValue "E29" '  -0.54%  '

Set-TextValue "D30" '4.420'
Set-TextValue "E30" '  +0.28%  '

Set-TextValue "D31" '4.342'
Set-TextValue "E31" '  +0.12%  '

Set-TextValue "D32" '1.212'
Set-TextValue "E32" '  -5.09%  '

Set-TextValue "D33" '0.05376'
Set-TextValue "E33" '  +2.31%  '

Set-TextValue "D34" '1.952'
Set-TextValue "E34" '  +1.06%  '

Set-TextValue "D35" '1.183'
Set-TextValue "E35" '  +0.53%  '

Set-TextValue "D36" '0.7494'
Set-TextValue "E36" '  +1.16%  '

Set-TextValue "D37" '2.700'
Set-TextValue "E37" '  +0.47%  '

Set-TextValue "D38" '0.01888'
Set-TextValue "E38" '  +1.20%  '

Set-TextValue "D39" '1.287.34'
Set-TextValue "E39" '  +9.63%  '

Set-TextValue "D40" '2.748'
Set-TextValue "E40" '  +0.68%  '

Set-TextValue "D41" '6.603'
Set-TextValue "E41" '  +3.45%  '

Set-TextValue "B42" 'TrustWalletToken'
Set-TextValue "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" '0.9036'
Set-TextValue "E42" '  +1.93%  '

Set-TextValue "B43" 'Quant'
Set-TextValue "C43" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D43" '111.59'
Set-TextValue "E43" '  +4.93%  '

Set-TextValue "B44" 'Aave'
Set-TextValue "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '73.67'
Set-TextValue "E44" '  +0.59%  '

Set-TextValue "E45" '  +0.11%  '

Set-TextValue "E46" '  +7.09%  '

Set-TextValue "D47" '2.030.69'
Set-TextValue "E47" '  +0.13%  '

Set-TextValue "E48" '  -0.19%  '

Set-TextValue "D49" '0.5218'
Set-TextValue "E49" '  +0.28%  '

Set-TextValue "D50" '9.506'
Set-TextValue "E50" '  +1.25%  '

Set-TextValue "D51" '0.4398'
Set-TextValue "E51" '  +2.16%  '
